$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Remove the "Meta description" paragraph that currently sits right
#    after the title paragraph ("Play Cold Spell Slot for Free ...").
# ---------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description*") {
        $p.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 2) Insert a new bold paragraph ("Play Cold Spell Slot for Free ...")
#    right before the final paragraph (the image-prompt paragraph).
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

# Create a fresh, empty paragraph right before the last one.
$lastPara.Range.InsertParagraphBefore() | Out-Null

# The freshly-created paragraph is now at the old last-paragraph index.
$targetPara = $d.Paragraphs.Item($count)
$insertionRange = $targetPara.Range
$insertionRange.Collapse(1)

$newParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Cold Spell Slot for Free – Novomatic Fantasy Theme</w:t></w:r></w:p>
</w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@

$insertionRange.InsertXML($newParaXml)

# ---------------------------------------------------------------------
# 3) Replace the text of the final (image-prompt) paragraph with the
#    meta-description copy, keeping its existing italic formatting.
# ---------------------------------------------------------------------
$oldText = "Please create a cartoon-style feature image for the online slot game ""Cold Spell"". The image should feature a happy Maya warrior with glasses. The Maya warrior should be depicted holding a wand and standing in front of ice-covered mountains with a snowy background to reflect the game's medieval fantasy theme. The image can include other elements from the game such as playing cards, tiaras, maps, and treasure chests. The image should be bright and colorful to capture the attention of players and entice them to try the game."
$newText = "Explore a medieval realm with Cold Spell, a Novomatic online slot game with stunning ice-covered reels and exciting win potential. Play for free now."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
